# developer guide: fixed typo and update diagrams
#
# 1) Bump every cached "datetimeFigureOut" date field (slide master, all
#    slide layouts, and the notes master) from 3/25/2019 to 3/31/2019.
# 2) Fix the "ProkectCard" -> "ProjectCard" typo on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "3/25/2019"
$newDate = "3/31/2019"
$ppPlaceholderDate = 16

function Update-DateShape($shp) {
    if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# Every slide layout off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# Notes master
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    Update-DateShape $notesMaster.Shapes.Item($i)
}

# Fix the "ProkectCard" typo on slide 1 (shape id 46, "Rectangle 11").
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "ProkectCard") {
            $shp.TextFrame.TextRange.Text = "ProjectCard"
        }
    }
}
